# The sheet tracks one price-report row per date for "Vega Modelo de Temuco - Coco".
# This commit adds a new weekly report row, inserted at row 15 (so the data stays in
# reverse-chronological-ish order as in the source feed), pushing the previously
# existing rows 15-117 down to 16-118.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 15; existing rows 15:117 shift down to 16:118.
$ws.Rows("15:15").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Vega Modelo de Temuco"
$ws.Range("C15").Value = "La Araucanía"
$ws.Range("D15").Value = 45149
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100108
$ws.Range("H15").Value = "Tropicales y subtropicales"
$ws.Range("I15").Value = 100108007
$ws.Range("J15").Value = "Coco"
$ws.Range("K15").Value = "Sin especificar"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 32000
$ws.Range("O15").Value = 32000
$ws.Range("P15").Value = 32000
$ws.Range("Q15").Value = "$/malla 20 unidades"
$ws.Range("R15").Value = "Perú"
$ws.Range("S15").Value = 1600
$ws.Range("T15").Value = 20
